# Daily attendance processing - 2025-11-15 11:18:51
# For every row in the "Recorded By" column (G), the last two comma-separated
# names/emails are swapped in order (their relative position in the list is
# flipped) while any preceding entries keep their place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Row + $used.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = 7
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $parts = $val.Split(",")
        $n = $parts.Length

        if ($n -ge 2) {
            for ($i = 0; $i -lt $n; $i++) {
                $parts[$i] = $parts[$i].Trim()
            }

            $tmp = $parts[$n - 1]
            $parts[$n - 1] = $parts[$n - 2]
            $parts[$n - 2] = $tmp

            $cell.Value2 = [string]::Join(", ", $parts)
        }
    }
}
